$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("7.745759138573715", "1690.630179487911", "224661.9279532027", "15706235.67871721", "1064189171.088317", "1393436492059.6", "3.408420886080534E+16")
    3  = @("3.029404614813376", "493.7056657904255", "51365.92510653412", "3946194.012936089", "301596277.797565", "491505184098.5649", "1.516427993093481E+16")
    4  = @("8.659901630721464", "1842.921612235855", "239041.5701604958", "16400704.41870035", "1100385677.789406", "1424477932023.208", "3.455941848899456E+16")
    5  = @("2.65747106496605", "1032.364311169534", "173959.1629817165", "14005904.0100794", "992793839.8827198", "1345866062803.361", "3.397062752496482E+16")
    6  = @("6.079569206452023", "623.1497992270617", "64392.06810269794", "4030683.788336271", "283645076.2654231", "433742937238.8953", "1.4192001549645E+16")
    7  = @("1.134104262201696", "58.16326180633573", "4680.645402168596", "401817.5202235829", "31959278.30984279", "53727782889.18716", "1738799778604794")
    8  = @("4.015035218607613", "218.4341029728908", "30721.6938354555", "2576969.00124671", "202282824.2453786", "334035903055.7594", "1.065821675356784E+16")
    9  = @("2.691116401345114", "606.379519325962", "94684.81007314504", "7647897.093405686", "573502194.401401", "869544076802.2577", "2.470795691006535E+16")
    10 = @("0.5495072954475011", "28.8466065137616", "27.25984047594418", "27.04997347839618", "34.51632402078002", "66863.28710318518", "1890293719.477777")
}

$columns = @("B", "C", "D", "E", "F", "G", "H")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $col = $columns[$i]
        $ws.Range("$col$row").Value = [double]$values[$i]
    }
}
